$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 132, shifting existing rows 132:154 down to 133:155
$ws.Rows("132:132").Insert()

# Populate the newly inserted row 132 with the new record
$ws.Range("A132").Value = 5
$ws.Range("B132").Value = "Macroferia Regional de Talca"
$ws.Range("C132").Value = "Maule"
$ws.Range("D132").Value = 45258
$ws.Range("E132").Value = 7
$ws.Range("F132").Value = 100112026
$ws.Range("G132").Value = "Haba"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 200
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = 12000
$ws.Range("N132").Value = "$/saco 25 kilos"
$ws.Range("O132").Value = "Región del Maule"
$ws.Range("P132").Value = 480
$ws.Range("Q132").Value = 25
$ws.Range("R132").Value = "Hortaliza"
